$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B23: reformat the existing "Finalisation de la revue 1. / Passage de la revue 1" text.
# Text stays the same overall, but the rich-text run split & font size change (11 -> 14),
# with the leading run losing its trailing space (moved into its own italic run).
$ws.Range("B23").Value = "Finalisation de la revue 1. Passage de la revue 1"

$part1 = "Finalisation de la revue 1."
$part2 = " "
$part3 = "Passage de la revue 1"
$start2 = $part1.Length + 1
$len2 = $part2.Length
$start3 = $start2 + $len2
$len3 = $part3.Length

$run2 = $ws.Range("B23").Characters($start2, $len2)
$run2.Font.Italic = $true
$run2.Font.Size = 14

$run3 = $ws.Range("B23").Characters($start3, $len3)
$run3.Font.Bold = $true
$run3.Font.Italic = $true
$run3.Font.Size = 14

# --- Row 45: the "Reprise du diaporama pour la revu de projet n°2." entry (H45/I45)
# moves away from this row; C45 becomes a blank, formatted (time) cell.
$ws.Range("H45").ClearContents()
$ws.Range("I45").Clear()
$ws.Range("C45").NumberFormat = $ws.Range("C44").NumberFormat
$ws.Range("C45").HorizontalAlignment = $ws.Range("C2").HorizontalAlignment
$ws.Range("C45").VerticalAlignment = $ws.Range("C2").VerticalAlignment

# --- Row 51: new B/C entries (student 1) + the relocated H/I entry (student 3).
$ws.Range("B51").Value = "recherche pour client/server socket"
$ws.Range("C51").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("C51").HorizontalAlignment = $ws.Range("C2").HorizontalAlignment
$ws.Range("C51").VerticalAlignment = $ws.Range("C2").VerticalAlignment
$ws.Range("C51").Value = 2/24

$ws.Range("H51").Value = "Reprise du diaporama pour la revu de projet n°2."
$ws.Range("I51").NumberFormat = $ws.Range("I9").NumberFormat
$ws.Range("I51").HorizontalAlignment = $ws.Range("I9").HorizontalAlignment
$ws.Range("I51").VerticalAlignment = $ws.Range("I9").VerticalAlignment
$ws.Range("I51").Value = 2/24

# --- Row 52: new B/C entries (student 1) + new I entry (student 3, no task text).
$ws.Range("B52").Value = "recherche pour client/server socket"
$ws.Range("C52").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("C52").HorizontalAlignment = $ws.Range("C2").HorizontalAlignment
$ws.Range("C52").VerticalAlignment = $ws.Range("C2").VerticalAlignment
$ws.Range("C52").Value = 3/24

$ws.Range("I52").NumberFormat = $ws.Range("I9").NumberFormat
$ws.Range("I52").HorizontalAlignment = $ws.Range("I9").HorizontalAlignment
$ws.Range("I52").VerticalAlignment = $ws.Range("I9").VerticalAlignment
$ws.Range("I52").Value = 3/24

# --- Row 54: new B/C entry (student 1, wrapped+centered style) + F54/I54 value bump.
$ws.Range("B54").Value = "Finalisation du digramme de squence pour le serveur socket asynchrone + préparation du diapo pour la revue 2 "
$ws.Range("B54").HorizontalAlignment = -4108
$ws.Range("B54").VerticalAlignment = -4108
$ws.Range("B54").WrapText = $true

$ws.Range("C54").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("C54").HorizontalAlignment = $ws.Range("C2").HorizontalAlignment
$ws.Range("C54").VerticalAlignment = $ws.Range("C2").VerticalAlignment
$ws.Range("C54").Value = 6/24

$ws.Range("F54").Value = 6/24

$ws.Range("I54").NumberFormat = $ws.Range("I9").NumberFormat
$ws.Range("I54").HorizontalAlignment = $ws.Range("I9").HorizontalAlignment
$ws.Range("I54").VerticalAlignment = $ws.Range("I9").VerticalAlignment
$ws.Range("I54").Value = 6/24

# --- Cursor/selection position as last saved by the author.
$ws.Range("B62").Select() | Out-Null
